$wb = $excel.ActiveWorkbook

# The "Single int" sheet's selection moves off the previous range once the
# new sheet becomes active.
$sheetSingleInt = $wb.Worksheets.Item("Single int")
$sheetSingleInt.Range("E3").Select()

# Add the new worksheet at the end of the workbook and name it "Sheet1"
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Sheet1"

# Header text above the table
$ws.Range("E1").Value = "Benchmarking sequence of different length:"

# Column headers (row 4)
$ws.Range("F4").Value = "VC8SP1"
$ws.Range("G4").Value = "gcc 4.4.0 (32)"
$ws.Range("H4").Value = "VC++ 10 (32)"
$ws.Range("I4").Value = "Intel 11.1 (32)"
$ws.Range("J4").Value = "gcc 4.4.0 (64)"
$ws.Range("K4").Value = "VC++ 10 (64)"
$ws.Range("L4").Value = "Intel 11.1 (64)"

# Sequence length column (E5:E12)
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 4
$ws.Range("E8").Value = 5
$ws.Range("E9").Value = 6
$ws.Range("E10").Value = 7
$ws.Range("E11").Value = 8
$ws.Range("E12").Value = 9

# VC8SP1 column (F5:F12)
$ws.Range("F5").Value = 0.256
$ws.Range("F6").Value = 0.7
$ws.Range("F7").Value = 1.044
$ws.Range("F8").Value = 1.443
$ws.Range("F9").Value = 1.76
$ws.Range("F10").Value = 2.214
$ws.Range("F11").Value = 2.756
$ws.Range("F12").Value = 3.29

# gcc 4.4.0 (32) column (G5:G12)
$ws.Range("G5").Value = 0.627
$ws.Range("G6").Value = 0.956
$ws.Range("G7").Value = 1.298
$ws.Range("G8").Value = 1.65
$ws.Range("G9").Value = 1.906
$ws.Range("G10").Value = 2.478
$ws.Range("G11").Value = 2.676
$ws.Range("G12").Value = 2.945

# VC++ 10 (32) column (H5:H12)
$ws.Range("H5").Value = 0.424
$ws.Range("H6").Value = 0.492
$ws.Range("H7").Value = 1.217
$ws.Range("H8").Value = 1.789
$ws.Range("H9").Value = 1.986
$ws.Range("H10").Value = 2.513
$ws.Range("H11").Value = 2.829
$ws.Range("H12").Value = 3.732

# Intel 11.1 (32) column (I5:I12)
$ws.Range("I5").Value = 0.569
$ws.Range("I6").Value = 0.98
$ws.Range("I7").Value = 1.448
$ws.Range("I8").Value = 1.989
$ws.Range("I9").Value = 2.596
$ws.Range("I10").Value = 3.242
$ws.Range("I11").Value = 3.559
$ws.Range("I12").Value = 4.246

# gcc 4.4.0 (64) column (J5:J12)
$ws.Range("J5").Value = 0.819
$ws.Range("J6").Value = 1.277
$ws.Range("J7").Value = 1.472
$ws.Range("J8").Value = 2.229
$ws.Range("J9").Value = 2.709
$ws.Range("J10").Value = 3.205
$ws.Range("J11").Value = 3.377
$ws.Range("J12").Value = 3.592

# VC++ 10 (64) column (K5:K12)
$ws.Range("K5").Value = 0.311
$ws.Range("K6").Value = 0.441
$ws.Range("K7").Value = 0.768
$ws.Range("K8").Value = 0.975
$ws.Range("K9").Value = 1.216
$ws.Range("K10").Value = 1.634
$ws.Range("K11").Value = 1.853
$ws.Range("K12").Value = 2.11

# Intel 11.1 (64) column (L5:L12)
$ws.Range("L5").Value = 0.426
$ws.Range("L6").Value = 0.695
$ws.Range("L7").Value = 1.019
$ws.Range("L8").Value = 1.265
$ws.Range("L9").Value = 1.634
$ws.Range("L10").Value = 2.008
$ws.Range("L11").Value = 2.324
$ws.Range("L12").Value = 2.662

$ws.Range("G5").Select()
